# Fixture edit: "adding test for fuzzy worksheet name matching"
#
# 1. Rename sheet "Root" -> "Main root" (all existing defined names that
#    pointed at Root!... follow the rename automatically).
# 2. Add a second, duplicate "_xlnm._FilterDatabase" defined name (suffixed
#    "_0") scoped to each of the four sheets, pointing at the same range as
#    the sheet's existing (non-hidden) _FilterDatabase name - this produces
#    the duplicate-primary-key fixture the commit message refers to.
# 3. Reset every sheet's zoom to 100% and make "Main root" the active /
#    selected tab (selection parked on B10), matching the recorded view
#    state.

$wb = $excel.ActiveWorkbook

function QuoteSheetName($name) {
    if ($name -match '[^A-Za-z0-9_]') {
        return "'" + $name + "'"
    }
    return $name
}

# --- 1. rename the "Root" sheet -----------------------------------------
$rootSheet = $wb.Worksheets.Item("Root")
$rootSheet.Name = "Main root"

# --- 2. duplicate the _FilterDatabase defined name on every sheet -------
$filterRanges = @{
    "Main root"        = "`$A`$1:`$B`$2"
    "Nodes"            = "`$A`$1:`$D`$4"
    "Leaves"           = "`$A`$1:`$F`$7"
    "One to many rows" = "`$A`$1:`$A`$13"
}

$sheetOrder = @("Main root", "Nodes", "Leaves", "One to many rows")
foreach ($sname in $sheetOrder) {
    $s = $wb.Worksheets.Item($sname)
    $addr = $filterRanges[$sname]
    $refersTo = "=" + (QuoteSheetName $sname) + "!" + $addr
    $s.Names.Add("_xlnm._FilterDatabase_0", $refersTo)
}

# --- 3. view state: zoom every sheet to 100%, re-home the active tab ---
foreach ($sname in @("Nodes", "Leaves", "One to many rows")) {
    $s = $wb.Worksheets.Item($sname)
    $s.Activate()
    $excel.ActiveWindow.Zoom = 100
}

$mainRoot = $wb.Worksheets.Item("Main root")
$mainRoot.Activate()
$excel.ActiveWindow.Zoom = 100
$mainRoot.Range("B10").Select() | Out-Null
